$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-10 Wednesday" "2025-12-11 Thursday"

Replace-Text "15×15=225" "74×78=5772"
Replace-Text "90×43=3870" "77×22=1694"
Replace-Text "24×99=2376" "42×46=1932"
Replace-Text "23×43=989" "27×96=2592"
Replace-Text "50×94=4700" "20×50=1000"

Replace-Text "24×23=552" "24×48=1152"
Replace-Text "77×59=4543" "72×82=5904"
Replace-Text "97×79=7663" "44×53=2332"
Replace-Text "86×33=2838" "41×35=1435"
Replace-Text "33×39=1287" "99×63=6237"

Replace-Text "48×41=1968" "31×71=2201"
Replace-Text "13×53=689" "26×16=416"
Replace-Text "96×53=5088" "24×19=456"
Replace-Text "14×77=1078" "31×38=1178"
Replace-Text "13×26=338" "15×74=1110"

Replace-Text "31×50=1550" "51×88=4488"
Replace-Text "36×16=576" "59×42=2478"
Replace-Text "85×54=4590" "80×92=7360"
Replace-Text "44×65=2860" "14×53=742"
Replace-Text "24×56=1344" "47×23=1081"

Replace-Text "88×11=968" "81×81=6561"
Replace-Text "43×63=2709" "41×41=1681"
Replace-Text "63×93=5859" "40×83=3320"
Replace-Text "91×33=3003" "95×73=6935"
Replace-Text "92×83=7636" "95×89=8455"
